{"js": "// \"Fig 3 p-values\" table: update the Hake column (3rd column) p-values.\n// Table layout: row 0 = header (\"\", \"Cod\", \"Hake\"); rows 1-7 = data rows, each\n// [label, Cod p-value, Hake p-value]. Only the Hake p-value (column index 2)\n// changes per row.\nconst table = context.document.body.tables.getFirst();\n\nconst hakeChanges = [\n  { row: 1, oldValue: \"0.06\", newValue: \"0.19\" },   // GDP 2016\n  { row: 2, oldValue: \"0.02\", newValue: \"<0.01\" },  // OHI 2016\n  { row: 3, oldValue: \"0.45\", newValue: \"0.89\" },   // OHI economic 2016\n  { row: 4, oldValue: \"0.11\", newValue: \"0.05\" },   // Technical Development\n  { row: 5, oldValue: \"0.87\", newValue: \"0.25\" },   // Compilance (scores)\n  { row: 6, oldValue: \"0.13\", newValue: \"0.12\" },   // Readiness\n  { row: 7, oldValue: \"0.02\", newValue: \"<0.01\" },  // Vulnerability\n];\n\nconst cells = hakeChanges.map((change) => table.getCell(change.row, 2));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nhakeChanges.forEach((change, i) => {\n  const cell = cells[i];\n  if (cell.value === change.oldValue) {\n    cell.value = change.newValue;\n  }\n});\nawait context.sync();\n", "ps1": "# \"Fig 3 p-values\" table: update the Hake column (3rd column) p-values.\n# Table layout: row 1 = header (\"\", \"Cod\", \"Hake\"); rows 2-8 = data rows, each\n# [label, Cod p-value, Hake p-value]. Only the Hake p-value changes per row:\n#   GDP 2016                -> 0.06  => 0.19\n#   OHI 2016                -> 0.02  => <0.01\n#   OHI economic 2016       -> 0.45  => 0.89\n#   Technical Development   -> 0.11  => 0.05\n#   Compilance (scores)     -> 0.87  => 0.25\n#   Readiness               -> 0.13  => 0.12\n#   Vulnerability           -> 0.02  => <0.01\n# Each \"old\" value below is unique in the document except \"0.02\", which\n# occurs in exactly the two rows (OHI 2016, Vulnerability) that both need\n# to become \"<0.01\" - so a whole-document, whole-word replace is safe and\n# unambiguous for every value.\n$d = $word.ActiveDocument\n\n$pValueChanges = @(\n    @{ Old = \"0.06\"; New = \"0.19\" },\n    @{ Old = \"0.02\"; New = \"<0.01\" },\n    @{ Old = \"0.45\"; New = \"0.89\" },\n    @{ Old = \"0.11\"; New = \"0.05\" },\n    @{ Old = \"0.87\"; New = \"0.25\" },\n    @{ Old = \"0.13\"; New = \"0.12\" }\n)\n\nforeach ($change in $pValueChanges) {\n    $range = $d.Content\n    $range.Find.Execute($change.Old, $false, $true, $false, $false, $false, $true, 1, $false, $change.New, 2)\n}\n"}
